$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.436534333333333
$ws.Range("H2").Value = 4.309603
$ws.Range("I2").Value = 0.03241561610838976
$ws.Range("J2").Value = 0.03241561610838976
$ws.Range("M2").Value = 18.43631966666667
$ws.Range("N2").Value = 55.308959
$ws.Range("O2").Value = 0.6034704469962782
$ws.Range("P2").Value = 0.6034704469962781
$ws.Range("Q2").Value = 26.48440618147522
$ws.Range("R2").Value = 238.359655633277
$ws.Range("S2").Value = 0.01956186634258972
$ws.Range("T2").Value = 0.01956186634258972
$ws.Range("G3").Value = 1.436534333333333
$ws.Range("H3").Value = 4.309603
$ws.Range("I3").Value = 0.03241561610838976
$ws.Range("J3").Value = 0.03241561610838976
$ws.Range("O3").Value = 0.1750419652256785
$ws.Range("P3").Value = 0.1750419652256784
$ws.Range("Q3").Value = 7.682037337396111
$ws.Range("R3").Value = 69.13833603656499
$ws.Range("S3").Value = 0.005674093147613702
$ws.Range("T3").Value = 0.005674093147613701
$ws.Range("G4").Value = 1.436534333333333
$ws.Range("H4").Value = 4.309603
$ws.Range("I4").Value = 0.03241561610838976
$ws.Range("J4").Value = 0.03241561610838976
$ws.Range("M4").Value = 6.766555
$ws.Range("N4").Value = 20.299665
$ws.Range("O4").Value = 0.2214875877780434
$ws.Range("P4").Value = 0.2214875877780434
$ws.Range("Q4").Value = 9.720388575888334
$ws.Range("R4").Value = 87.48349718299501
$ws.Range("S4").Value = 0.007179656618186333
$ws.Range("T4").Value = 0.007179656618186332
$ws.Range("G5").Value = 34.88211266666666
$ws.Range("I5").Value = 0.7871201871162607
$ws.Range("J5").Value = 0.7871201871162609
$ws.Range("M5").Value = 18.43631966666667
$ws.Range("N5").Value = 55.308959
$ws.Range("O5").Value = 0.6034704469962782
$ws.Range("P5").Value = 0.6034704469962781
$ws.Range("Q5").Value = 643.097779771349
$ws.Range("R5").Value = 5787.880017942142
$ws.Range("S5").Value = 0.475003771158844
$ws.Range("T5").Value = 0.475003771158844
$ws.Range("G6").Value = 34.88211266666666
$ws.Range("I6").Value = 0.7871201871162607
$ws.Range("J6").Value = 0.7871201871162609
$ws.Range("O6").Value = 0.1750419652256785
$ws.Range("P6").Value = 0.1750419652256784
$ws.Range("Q6").Value = 186.5362252016655
$ws.Range("S6").Value = 0.137779064421634
$ws.Range("T6").Value = 0.137779064421634
$ws.Range("G7").Value = 34.88211266666666
$ws.Range("I7").Value = 0.7871201871162607
$ws.Range("J7").Value = 0.7871201871162609
$ws.Range("M7").Value = 6.766555
$ws.Range("N7").Value = 20.299665
$ws.Range("O7").Value = 0.2214875877780434
$ws.Range("P7").Value = 0.2214875877780434
$ws.Range("Q7").Value = 236.0317338751967
$ws.Range("R7").Value = 2124.28560487677
$ws.Range("S7").Value = 0.1743373515357828
$ws.Range("T7").Value = 0.1743373515357828
$ws.Range("G8").Value = 7.997472999999999
$ws.Range("H8").Value = 23.992419
$ws.Range("I8").Value = 0.1804641967753495
$ws.Range("J8").Value = 0.1804641967753495
$ws.Range("M8").Value = 18.43631966666667
$ws.Range("N8").Value = 55.308959
$ws.Range("O8").Value = 0.6034704469962782
$ws.Range("P8").Value = 0.6034704469962781
$ws.Range("Q8").Value = 147.4439687535356
$ws.Range("R8").Value = 1326.995718781821
$ws.Range("S8").Value = 0.1089048094948444
$ws.Range("T8").Value = 0.1089048094948444
$ws.Range("G9").Value = 7.997472999999999
$ws.Range("H9").Value = 23.992419
$ws.Range("I9").Value = 0.1804641967753495
$ws.Range("J9").Value = 0.1804641967753495
$ws.Range("O9").Value = 0.1750419652256785
$ws.Range("P9").Value = 0.1750419652256784
$ws.Range("Q9").Value = 42.76743323513833
$ws.Range("R9").Value = 384.906899116245
$ws.Range("S9").Value = 0.03158880765643072
$ws.Range("T9").Value = 0.03158880765643071
$ws.Range("G10").Value = 7.997472999999999
$ws.Range("H10").Value = 23.992419
$ws.Range("I10").Value = 0.1804641967753495
$ws.Range("J10").Value = 0.1804641967753495
$ws.Range("M10").Value = 6.766555
$ws.Range("N10").Value = 20.299665
$ws.Range("O10").Value = 0.2214875877780434
$ws.Range("P10").Value = 0.2214875877780434
$ws.Range("Q10").Value = 54.115340915515
$ws.Range("R10").Value = 487.038068239635
$ws.Range("S10").Value = 0.03997057962407431
$ws.Range("T10").Value = 0.03997057962407431

Write-Output "Applied 98 cell updates"
